# "perskirstytos valandos" - redistribute planned hours across tasks
# on the "Week" sheet's Development Tasks Completed table (Plan Hours column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week")

# Userių paskyrų trinimas: Plan Hours 3 -> 2
$ws.Range("C27").Value = 2

# Userių promotion/demotion: Plan Hours 4 -> 2
$ws.Range("C28").Value = 2

# Prekių peržiūra, dėjimas į krepšelį: Plan Hours 4 -> 5
$ws.Range("C29").Value = 5

# Testavimas: Plan Hours 3 -> 5
$ws.Range("C31").Value = 5

# Reflect the author's final cursor position on the sheet
$ws.Range("C29").Select() | Out-Null
